$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 89.64706
$ws.Range("I9").Value = 90.166664
$ws.Range("J9").Value = 88.40000000000001
$ws.Range("K9").Value = 90.166664
$ws.Range("L9").Value = 88.40000000000001
$ws.Range("M9").Value = 78.833336
$ws.Range("N9").Value = -426.4

$ws.Range("H28").Value = 459.70587
$ws.Range("I28").Value = 365.35715
$ws.Range("J28").Value = 900
$ws.Range("K28").Value = 365.35715
$ws.Range("L28").Value = 900
$ws.Range("M28").Value = 119.64285
$ws.Range("N28").Value = -1870

$ws.Range("H76").Value = 7121.552
$ws.Range("I76").Value = 9363.8125
$ws.Range("J76").Value = 4361.846
$ws.Range("K76").Value = 9363.8125
$ws.Range("L76").Value = 4361.846
$ws.Range("M76").Value = -9048.8125
$ws.Range("N76").Value = -4991.846

$ws.Range("H79").Value = 7121.552
$ws.Range("I79").Value = 9363.8125
$ws.Range("J79").Value = 4361.846
$ws.Range("K79").Value = 9363.8125
$ws.Range("L79").Value = 4361.846
$ws.Range("M79").Value = -8271.8125
$ws.Range("N79").Value = -6545.846

$ws.Range("H86").Value = 7534.5
$ws.Range("I86").Value = 8549.9375
$ws.Range("J86").Value = 4826.6665
$ws.Range("K86").Value = 8549.9375
$ws.Range("L86").Value = 4826.6665
$ws.Range("M86").Value = -7426.9375
$ws.Range("N86").Value = -7072.6665

$ws.Range("H89").Value = 7534.5
$ws.Range("I89").Value = 8549.9375
$ws.Range("J89").Value = 4826.6665
$ws.Range("K89").Value = 42749.6875
$ws.Range("L89").Value = 24133.3325
$ws.Range("M89").Value = -37133.6875
$ws.Range("N89").Value = -35365.3325

$ws.Range("H112").Value = 2493.5483
$ws.Range("I112").Value = 400
$ws.Range("J112").Value = 3104.1667
$ws.Range("K112").Value = 1200
$ws.Range("L112").Value = 9312.500100000001
$ws.Range("M112").Value = -92
$ws.Range("N112").Value = -11528.5001

$ws.Range("H137").Value = 2263.8333
$ws.Range("I137").Value = 1563.6666
$ws.Range("J137").Value = 2497.2222
$ws.Range("K137").Value = 4690.9998
$ws.Range("L137").Value = 7491.6666
$ws.Range("M137").Value = -2140.9998
$ws.Range("N137").Value = -12591.6666

$ws.Range("H138").Value = 2370.0908
$ws.Range("I138").Value = 2204.625
$ws.Range("J138").Value = 2413.4917
$ws.Range("K138").Value = 6613.875
$ws.Range("L138").Value = 7240.4751
$ws.Range("M138").Value = -1473.875
$ws.Range("N138").Value = -17520.4751

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20100.934
$ws.Range("I32").Value = 3881.5818
$ws.Range("J32").Value = 64704.15
$ws.Range("K32").Value = 3881.5818
$ws.Range("L32").Value = 64704.15
$ws.Range("M32").Value = -3594.5818
$ws.Range("N32").Value = -65278.15

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2841.5334
$ws.Range("I31").Value = 1599.0278
$ws.Range("J31").Value = 7811.5557
$ws.Range("K31").Value = 1599.0278
$ws.Range("L31").Value = 7811.5557
$ws.Range("M31").Value = -1304.0278
$ws.Range("N31").Value = -8401.555700000001

$ws.Range("H34").Value = 2841.5334
$ws.Range("I34").Value = 1599.0278
$ws.Range("J34").Value = 7811.5557
$ws.Range("K34").Value = 1599.0278
$ws.Range("L34").Value = 7811.5557
$ws.Range("M34").Value = -1397.0278
$ws.Range("N34").Value = -8215.555700000001

$ws.Range("H122").Value = 2795.625
$ws.Range("I122").Value = 1617.091
$ws.Range("J122").Value = 3792.8462
$ws.Range("K122").Value = 4851.272999999999
$ws.Range("L122").Value = 11378.5386
$ws.Range("M122").Value = -2401.272999999999
$ws.Range("N122").Value = -16278.5386

$ws.Range("H132").Value = 4393.5
$ws.Range("I132").Value = 4080.4443
$ws.Range("J132").Value = 5332.6665
$ws.Range("K132").Value = 12241.3329
$ws.Range("L132").Value = 15997.9995
$ws.Range("M132").Value = -9711.332900000001
$ws.Range("N132").Value = -21057.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 853.9231
$ws.Range("I86").Value = 595.8570999999999
$ws.Range("J86").Value = 1155
$ws.Range("K86").Value = 1787.5713
$ws.Range("L86").Value = 3465
$ws.Range("M86").Value = -601.5712999999998
$ws.Range("N86").Value = -5837

$ws.Range("H87").Value = 11675
$ws.Range("I87").Value = 8350
$ws.Range("J87").Value = 15000
$ws.Range("K87").Value = 25050
$ws.Range("L87").Value = 45000
$ws.Range("M87").Value = -23802
$ws.Range("N87").Value = -47496

$ws.Range("H89").Value = 853.9231
$ws.Range("I89").Value = 595.8570999999999
$ws.Range("J89").Value = 1155
$ws.Range("K89").Value = 5362.7139
$ws.Range("L89").Value = 10395
$ws.Range("M89").Value = 565.2861000000003
$ws.Range("N89").Value = -22251

$ws.Range("H90").Value = 11675
$ws.Range("I90").Value = 8350
$ws.Range("J90").Value = 15000
$ws.Range("K90").Value = 75150
$ws.Range("L90").Value = 135000
$ws.Range("M90").Value = -68910
$ws.Range("N90").Value = -147480

$ws.Range("H92").Value = 950
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 950
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 2850
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -5346

$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()

$ws.Range("H94").Value = 6717.7144
$ws.Range("I94").Value = 4008
$ws.Range("J94").Value = 8750
$ws.Range("K94").Value = 12024
$ws.Range("L94").Value = 26250
$ws.Range("M94").Value = -11348
$ws.Range("N94").Value = -27602

$ws.Range("H95").Value = 11345
$ws.Range("I95").Value = 2690
$ws.Range("J95").Value = 20000
$ws.Range("K95").Value = 8070
$ws.Range("L95").Value = 60000
$ws.Range("M95").Value = -6011
$ws.Range("N95").Value = -64118

$ws.Range("H96").Value = 3896.1538
$ws.Range("I96").Value = 3000
$ws.Range("J96").Value = 3970.8333
$ws.Range("K96").Value = 9000
$ws.Range("L96").Value = 11912.4999
$ws.Range("M96").Value = -6941
$ws.Range("N96").Value = -16030.4999

$ws.Range("H97").Value = 451.85715
$ws.Range("I97").Value = 267.66666
$ws.Range("J97").Value = 590
$ws.Range("K97").Value = 802.9999799999999
$ws.Range("L97").Value = 1770
$ws.Range("M97").Value = -306.9999799999999
$ws.Range("N97").Value = -2762

$ws.Range("H98").Value = 276.63635
$ws.Range("I98").Value = 257.5
$ws.Range("J98").Value = 299.6
$ws.Range("K98").Value = 772.5
$ws.Range("L98").Value = 898.8000000000001
$ws.Range("M98").Value = 725.5
$ws.Range("N98").Value = -3894.8

$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()

$ws.Range("H101").Value = 10400
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 10400
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 31200
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -36068

$ws.Range("H102").Value = 6869.5
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 6869.5
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 20608.5
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -25476.5

$ws.Range("H122").Value = 3707.6924
$ws.Range("I122").Value = 500
$ws.Range("J122").Value = 3836
$ws.Range("K122").Value = 4500
$ws.Range("L122").Value = 34524
$ws.Range("M122").Value = -2050
$ws.Range("N122").Value = -39424

$ws.Range("H130").Value = 2000
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 2000
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 6000
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -16040

$ws.Range("H131").Value = 1038.7142
$ws.Range("I131").Value = 506
$ws.Range("J131").Value = 1127.5
$ws.Range("K131").Value = 1518
$ws.Range("L131").Value = 3382.5
$ws.Range("M131").Value = 3522
$ws.Range("N131").Value = -13462.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 706.62067
$ws.Range("I55").Value = 387.08334
$ws.Range("J55").Value = 932.17645
$ws.Range("K55").Value = 387.08334
$ws.Range("L55").Value = 932.17645
$ws.Range("M55").Value = -214.08334
$ws.Range("N55").Value = -1278.17645

$ws.Range("H94").Value = 25330
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 25330
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 25330
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -26682

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1269.6842
$ws.Range("I136").Value = 1033.8148
$ws.Range("J136").Value = 1848.6364
$ws.Range("K136").Value = 3101.4444
$ws.Range("L136").Value = 5545.9092
$ws.Range("M136").Value = -551.4444000000003
$ws.Range("N136").Value = -10645.9092
